$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.06493566666666667
$ws.Range("H2").Value = 0.194807
$ws.Range("I2").Value = 0.05459895593058446
$ws.Range("J2").Value = 0.05459895593058446
$ws.Range("Q2").Value = 0.01783310897488889
$ws.Range("R2").Value = 0.160497980774
$ws.Range("S2").Value = 0.05459895593058446
$ws.Range("T2").Value = 0.05459895593058446

# Row 3 updates
$ws.Range("I3").Value = 0.8534208043695533
$ws.Range("J3").Value = 0.8534208043695533
$ws.Range("S3").Value = 0.8534208043695533
$ws.Range("T3").Value = 0.8534208043695533

# Row 4 updates
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.109394
$ws.Range("H4").Value = 0.328182
$ws.Range("I4").Value = 0.09198023969986227
$ws.Range("J4").Value = 0.09198023969986227
$ws.Range("Q4").Value = 0.03004258250266666
$ws.Range("R4").Value = 0.270383242524
$ws.Range("S4").Value = 0.09198023969986227
$ws.Range("T4").Value = 0.09198023969986227
